$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the casing of the existing "ampl. Min" label before we shuffle rows
# around (so the shared-string table collapses the old entry and appends
# the corrected one, matching the target workbook).
$ws.Range("A7").Value = "ampl. min"

# Insert the new rows bottom-up so earlier row numbers stay stable while we work.

# New "ampl. max" row, right after the (soon to be renamed) "ampl. min" row,
# before the "fraction" row.
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "ampl. max"

# New "sigma max" row, right after "sigma min", before "amplitude".
$ws.Rows("6:6").Insert()
$ws.Range("A6").Value = "sigma max"

# New "center min" / "center max" rows, right after "center", before "sigma".
$ws.Rows("4:5").Insert()
$ws.Range("A4").Value = "center min"
$ws.Range("A5").Value = "center max"

# Match the saved selection/active cell from the authored workbook.
$ws.Range("B15").Select()
